# This workbook tracks daily "Zapallo italiano" price records for the
# Vega Central Mapocho de Santiago market. The commit adds two new daily
# price records, inserted as rows 346 and 347 (pushing the existing
# records that used to be rows 346-439 down to rows 348-441).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 346, shifting every
# row from 346 downward by two positions.
$ws.Rows("346:347").Insert()

# --- New record: row 346 ---
$ws.Range("A346").Value = 9
$ws.Range("B346").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C346").Value = "Metropolitana"
$ws.Range("D346").Value = 44855
$ws.Range("E346").Value = 13
$ws.Range("F346").Value = 100112032
$ws.Range("G346").Value = "Zapallo italiano"
$ws.Range("H346").Value = "Sin especificar"
$ws.Range("I346").Value = "Primera"
$ws.Range("J346").Value = 250
$ws.Range("K346").Value = 12000
$ws.Range("L346").Value = 13000
$ws.Range("M346").Value = 12400
$ws.Range("N346").Value = "`$/caja 50 unidades"
$ws.Range("O346").Value = "Región de O'Higgins"
$ws.Range("P346").Value = 248
$ws.Range("Q346").Value = 50
$ws.Range("R346").Value = "Hortaliza"

# --- New record: row 347 ---
$ws.Range("A347").Value = 9
$ws.Range("B347").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C347").Value = "Metropolitana"
$ws.Range("D347").Value = 44855
$ws.Range("E347").Value = 13
$ws.Range("F347").Value = 100112032
$ws.Range("G347").Value = "Zapallo italiano"
$ws.Range("H347").Value = "Sin especificar"
$ws.Range("I347").Value = "Primera"
$ws.Range("J347").Value = 130
$ws.Range("K347").Value = 15000
$ws.Range("L347").Value = 15000
$ws.Range("M347").Value = 15000
$ws.Range("N347").Value = "`$/caja 60 unidades"
$ws.Range("O347").Value = "Limache"
$ws.Range("P347").Value = 250
$ws.Range("Q347").Value = 60
$ws.Range("R347").Value = "Hortaliza"
